$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14, shifting rows 14-21 down to 15-22.
# (Topic 3 / NumPy tutorial gets its own spacer row, same as the other topics.)
$ws.Rows.Item(14).Insert()

# --- New row 23: "find a tool to extract documentation..." task ---
$ws.Range("A23").Value = 9
$ws.Range("B22").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "find a tool to extarct documentation from the docstrings to test it works ok"
$ws.Range("C22").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 43636
$ws.Range("F23").Value = "Ongoing "
$ws.Range("G23").Value = "use python help() function, as described here https://wiki.python.org/moin/DocumentationTools"
$ws.Rows.Item(23).RowHeight = $ws.Rows.Item(22).RowHeight

# Row 13: fill in the previously-blank "Topic 3" row
$ws.Range("C13").Value = 43636
$ws.Range("D13").Value = "Topic 3: NumPy tutorial"

# Row 17 (the "TestCases for the RigidBody Class" row): add new comment
$ws.Range("G17").Value = "have timing left "

# Row 22 (the "Add quaternions function" row, was row 21 before the insert):
# status ToDo -> Done, add comment
$ws.Range("F22").Value = "Done"
$ws.Range("G22").Value = " to investigate more numerically stable formulas, see if I can do better in tests"

# Row 21 (the "Improve the logger" row, was row 20 before the insert):
# status ToDo -> Done, add comment
$ws.Range("F21").Value = "Done"
$ws.Range("G21").Value = "enough for this time, open new issue again if needed"

# Update the sheet view / selection to match the final authored state
$ws.Range("B15").Select()
